# Update "北京-漫展信息" workbook to the newer scrape snapshot.
# Commit message: "Update gh-pages to output generated at 456a3b4"
#
# Summary of the change:
#  1. Sheet "展览" (Exhibitions)      - several "want to go" counts (col F) bumped up.
#  2. Sheet "演出" (Performances)     - two "want to go" counts (col F) bumped up.
#  3. Sheet "本地生活" (Local life)   - the oldest entry (剑网3 cafe) expired and
#     dropped off the feed, so every remaining row shifted up by one and picked
#     up refreshed data (dates / counts / links); the sheet shrank from 6 to 5
#     used rows.
#  4. Sheet "全部类型" (All types)    - mirror of the same "want to go" bumps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. 展览 (Exhibitions) - col F updates
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(4, 6).Value = 5742
$ws1.Cells.Item(6, 6).Value = 58
$ws1.Cells.Item(7, 6).Value = 535
$ws1.Cells.Item(9, 6).Value = 1553
$ws1.Cells.Item(10, 6).Value = 13
$ws1.Cells.Item(13, 6).Value = 1557
$ws1.Cells.Item(14, 6).Value = 1557
$ws1.Cells.Item(15, 6).Value = 1453
$ws1.Cells.Item(16, 6).Value = 537
$ws1.Cells.Item(18, 6).Value = 583
$ws1.Cells.Item(19, 6).Value = 4148
$ws1.Cells.Item(20, 6).Value = 4148
$ws1.Cells.Item(22, 6).Value = 3315
$ws1.Cells.Item(23, 6).Value = 787
$ws1.Cells.Item(25, 6).Value = 2250
$ws1.Cells.Item(27, 6).Value = 323
$ws1.Cells.Item(29, 6).Value = 440
$ws1.Cells.Item(30, 6).Value = 1201
$ws1.Cells.Item(33, 6).Value = 1117
$ws1.Cells.Item(34, 6).Value = 1134

# ---------------------------------------------------------------------------
# 2. 演出 (Performances) - col F updates
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(18, 6).Value = 279
$ws2.Cells.Item(19, 6).Value = 215

# ---------------------------------------------------------------------------
# 3. 本地生活 (Local life) - oldest row (剑网3 cafe) expired; remove it and
#    shift everything else up, then refresh the data that changed underneath.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

# Row 3 (剑网3 × HAPPY ZOO) has expired/rolled off - delete it, shifting the
# rows below (old 4,5,6) up into (3,4,5).
$ws3.Rows.Item(3).Delete()

# Deleting the row also shifted the plain numeric index in col A along with it
# (old A4=3 -> A3=3, etc); put the sequential index back (2,3,4) same as before.
$ws3.Cells.Item(3, 1).Value = 2
$ws3.Cells.Item(4, 1).Value = 3
$ws3.Cells.Item(5, 1).Value = 4

# Row 3 (was old row 4: EVANGELION x PrismLand) - refreshed counts/link.
# (leading apostrophe keeps the yyyy-mm-dd column stored as text, as it was
# originally, instead of being auto-converted to a date value)
$ws3.Cells.Item(3, 2).Value = "'2024-10-15"
$ws3.Cells.Item(3, 3).Value = "北京·EVANGELION× PrismLand · 新世纪福音战士官方授权主题店"
$ws3.Cells.Item(3, 4).Value = "王府井地铁站F1东口步行120米 北京王府井喜悦购物中心"
$ws3.Cells.Item(3, 5).Value = "2024.10.15 00:00-12.15 23:59"
$ws3.Cells.Item(3, 6).Value = 574
$ws3.Cells.Item(3, 7).Value = 20
$ws3.Cells.Item(3, 8).Value = "https://show.bilibili.com/platform/detail.html?id=93039"
$ws3.Cells.Item(3, 9).Value = "//i0.hdslb.com/bfs/openplatform/202409/n32CfRya1727584778969.jpeg"

# Row 4 (was old row 5: 全职高手 x HAPPY ZOO) - refreshed counts/link.
$ws3.Cells.Item(4, 2).Value = "'2024-10-17"
$ws3.Cells.Item(4, 3).Value = "北京·全职高手×HAPPY ZOO 全职高手十周年咖啡厅"
$ws3.Cells.Item(4, 4).Value = "学清路38号金码大厦B座(六道口地铁站B东北口步行110米) BOM嘻番里"
$ws3.Cells.Item(4, 5).Value = "2024.10.17 00:00-11.17 23:59"
$ws3.Cells.Item(4, 6).Value = 115
$ws3.Cells.Item(4, 7).Value = 10
$ws3.Cells.Item(4, 8).Value = "https://show.bilibili.com/platform/detail.html?id=93324"
$ws3.Cells.Item(4, 9).Value = "//i2.hdslb.com/bfs/openplatform/202410/bVeTwUWu1728699225130.png"

# Row 5 (was old row 6: 蜡笔小新 x HAPPY ZOO) - refreshed counts/link, now sold out.
$ws3.Cells.Item(5, 2).Value = "'2024-10-17"
$ws3.Cells.Item(5, 3).Value = "北京·蜡笔小新：我们的恐龙日记x HAPPY ZOO 主题咖啡厅"
$ws3.Cells.Item(5, 4).Value = "王府井地铁站F1东口步行120米 北京王府井喜悦购物中心"
$ws3.Cells.Item(5, 5).Value = "2024.10.17 00:00-10.27 23:59"
$ws3.Cells.Item(5, 6).Value = 207
$ws3.Cells.Item(5, 7).Value = "已售罄"
$ws3.Cells.Item(5, 8).Value = "https://show.bilibili.com/platform/detail.html?id=93224"
$ws3.Cells.Item(5, 9).Value = "//i0.hdslb.com/bfs/openplatform/202410/nzGP5KRA1728526131597.png"

# ---------------------------------------------------------------------------
# 4. 全部类型 (All types) - col F updates (mirrors the sheets above)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(6, 6).Value = 574
$ws4.Cells.Item(7, 6).Value = 115
$ws4.Cells.Item(8, 6).Value = 5742
$ws4.Cells.Item(11, 6).Value = 58
$ws4.Cells.Item(16, 6).Value = 535
$ws4.Cells.Item(19, 6).Value = 1553
$ws4.Cells.Item(21, 6).Value = 13
$ws4.Cells.Item(23, 6).Value = 1557
$ws4.Cells.Item(24, 6).Value = 1557
$ws4.Cells.Item(26, 6).Value = 1453
$ws4.Cells.Item(27, 6).Value = 537
$ws4.Cells.Item(29, 6).Value = 583
$ws4.Cells.Item(31, 6).Value = 4148
$ws4.Cells.Item(32, 6).Value = 4148
$ws4.Cells.Item(34, 6).Value = 3315
$ws4.Cells.Item(35, 6).Value = 787
$ws4.Cells.Item(37, 6).Value = 2250
$ws4.Cells.Item(40, 6).Value = 440
$ws4.Cells.Item(41, 6).Value = 1201
$ws4.Cells.Item(44, 6).Value = 279
$ws4.Cells.Item(45, 6).Value = 215
$ws4.Cells.Item(49, 6).Value = 1117
$ws4.Cells.Item(50, 6).Value = 1134
